# Draft results rewrite, FHY, input pyrite data
# Update the "Counts (cps)" (column B) and "Error (cps)" (column C) values
# in the "Equilibrated Data" sheet rows 2-19 with the newly re-processed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Equilibrated Data")

$ws.Range("B2").Value = 1.24833333333333
$ws.Range("C2").Value = 0.0912531666666667

$ws.Range("B3").Value = 1.217
$ws.Range("C3").Value = 0.09010668

$ws.Range("B4").Value = 1.25733333333333
$ws.Range("C4").Value = 0.09158416

$ws.Range("B5").Value = 4.271
$ws.Range("C5").Value = 0.16878992

$ws.Range("B6").Value = 4.36466666666667
$ws.Range("C6").Value = 0.170571173333333

$ws.Range("B7").Value = 4.63233333333333
$ws.Range("C7").Value = 0.175750726666667

$ws.Range("B8").Value = 17.603
$ws.Range("C8").Value = 0.34290644

$ws.Range("B9").Value = 17.7466666666667
$ws.Range("C9").Value = 0.3439304

$ws.Range("B10").Value = 17.9223333333333
$ws.Range("C10").Value = 0.345542586666667

$ws.Range("B11").Value = 33.962
$ws.Range("C11").Value = 0.47614724

$ws.Range("B12").Value = 34.8266666666667
$ws.Range("C12").Value = 0.482001066666667

$ws.Range("B13").Value = 34.5783333333333
$ws.Range("C13").Value = 0.480638833333333

$ws.Range("B14").Value = 78.6723333333333
$ws.Range("C14").Value = 0.725358913333333

$ws.Range("B15").Value = 86.241
$ws.Range("C15").Value = 0.75719598

$ws.Range("B16").Value = 86.4766666666667
$ws.Range("C16").Value = 0.760994666666667

$ws.Range("B17").Value = 168.513333333333
$ws.Range("C17").Value = 1.061634

$ws.Range("B18").Value = 166.426333333333
$ws.Range("C18").Value = 1.05514295333333

$ws.Range("B19").Value = 166.140333333333
$ws.Range("C19").Value = 1.05000690666667
